$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the bold/bordered header-row formatting (A1:K1 reverts to the
#     default "Normal" style, matching the rest of the sheet) ---
$ws.Range("A1:K1").ClearFormats()

# --- Append four new data rows (38-41) below the existing data, mirroring
#     the pattern used by the other rows in the sheet ---
$ws.Range("A38").Value = "'08.14.18"
$ws.Range("B38").Value = "H.BROWN"
$ws.Range("C38").Value = 38
$ws.Range("D38").Value = 3275
$ws.Range("E38").Value = "fullRNASeq"
$ws.Range("F38").Value = "Brent_3275_17-5_GTAC_5_NEB_Universal_ATCGAGC_AGATCTCG_S6_R1_001"
$ws.Range("A38").ClearFormats()

$ws.Range("A39").Value = "'10.18.18"
$ws.Range("B39").Value = "H.BROWN"
$ws.Range("C39").Value = 39
$ws.Range("D39").Value = 3275
$ws.Range("E39").Value = "fullRNASeq"
$ws.Range("F39").Value = "Brent_3275_18-5_GTAC_5_SIC_Index2_010_ATCGAGC_GCTTCTAG_S15_R1_001"
$ws.Range("A39").ClearFormats()

$ws.Range("A40").Value = "'11.02.18"
$ws.Range("B40").Value = "H.BROWN"
$ws.Range("C40").Value = 40
$ws.Range("D40").Value = 3275
$ws.Range("E40").Value = "fullRNASeq"
$ws.Range("F40").Value = "Brent_3275_20-4_GTAC_4_SIC_Index2_07_CACCTCC_GAGTTGAG_S34_R1_001"
$ws.Range("A40").ClearFormats()

$ws.Range("A41").Value = "'10.18.18"
$ws.Range("B41").Value = "H.BROWN"
$ws.Range("C41").Value = 41
$ws.Range("D41").Value = 3275
$ws.Range("E41").Value = "fullRNASeq"
$ws.Range("F41").Value = "Brent_3275_18-3_GTAC_3_SIC_Index2_010_ATGACAG_GCTTCTAG_S13_R1_001"
$ws.Range("A41").ClearFormats()

# --- Selection cursor moves to O15 (matches the saved view state) ---
[void]$ws.Range("O15").Select()
